$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Valid" -> "Result"
$ws.Range("C1").Value = "Result"

# Row 2: John Doe -> Oogway Doe, and result Yes -> Valid
$ws.Range("A2").Value = "OogwayDoe@gmail.com"
$ws.Range("B2").Value = "Oogway@123"
$ws.Range("C2").Value = "Valid"

# Row 3: result No -> Invalid
$ws.Range("C3").Value = "Invalid"

# Row 4: result No -> Invalid
$ws.Range("C4").Value = "Invalid"
